# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    matching the header layout used by the other quarterly sheets, and
#    populate it with the Q1-2022 fund holdings.
# 2. Insert a new top row into "总计" for the "2022-Q1" aggregate and
#    renumber the existing rows' index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet (positioned immediately before "总计")
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Copy the header-row + index-column formatting from the "2021-Q4" sheet,
# which uses the same column layout (基金代码/基金名称/基金规模/...).
$fmtSrc = $wb.Worksheets.Item("2021-Q4")
$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows. 基金代码 (B) and 基金规模/股票总仓位/仓位占比/持有市值 (D:G) are
# stored as text (as in every other quarter's sheet), A (index) and H
# (仓位排名) are stored as numbers.
$rows = @(
    @(0, "005176", "富国精准医疗灵活配置混合", "33.81", "93.54", "3.70", "1.2510", 8),
    @(1, "000452", "南方医药保健灵活配置混合", "31.88", "90.98", "3.39", "1.0807", 10),
    @(2, "010703", "财通智选消费股票A", "1.12", "92.43", "3.57", "0.0400", 2),
    @(3, "010704", "财通智选消费股票C", "0.54", "92.43", "3.57", "0.0193", 2),
    @(4, "008884", "博远博锐混合A", "0.19", "86.59", "4.76", "0.0090", 5),
    @(5, "004917", "中银证券祥瑞混合A", "0.10", "79.01", "2.71", "0.0027", 3),
    @(6, "004918", "中银证券祥瑞混合C", "0.07", "79.01", "2.71", "0.0019", 3),
    @(7, "008885", "博远博锐混合C", "0.02", "86.59", "4.76", "0.0010", 5)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]

    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[1]

    $q1.Cells.Item($r, 3).Value = $row[2]

    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row[3]

    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row[4]

    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row[5]

    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row[6]

    $q1.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend the "2022-Q1" summary row to "总计" and renumber the index
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$tot.Range("A2:D2").Insert(-4121)
$tot.Range("B2:D2").ClearFormats()

$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 8
$tot.Range("D2").Value = 2.41

$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3

# Keep the originally-active sheet/tab selection untouched.
$wb.Worksheets.Item("2021-Q2").Activate()
